# Add a new git_lib entry on building a git server (row 48 on Sheet1),
# matching the three new shared-string rows and the new table row from
# the diff, then move the sheet's scroll/selection state to reflect the
# freshly-added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New data row -----------------------------------------------------
# Column A: "Language" value, Column B: short "Desc", Column C: the
# long markdown "Code Ref" body (wrapped, like the rows above it).
$ws.Range("A48").Value = "Remote "
$ws.Range("B48").Value = "My own git server"

$codeRef = @'
# Different Ways to Host Git Server
* Local Protocol: Local File System / Network File System
* SSH Protocol
* Http Protocol
* Git Protocol
# Use Local Protocol
* Clone existing project to build: `git clone --bare my_project my_project.git`
* Init an Empty project git: `git init --bare new_project.git`
'@
$ws.Range("C48").Value = $codeRef

# Match the wrapped-text style already used by column C on the rows
# above (style index 7 in the original workbook).
$ws.Range("C48").WrapText = $true

# Row height for the new row (matches the taller markdown body).
$ws.Rows.Item(48).RowHeight = 120

# --- View state ---------------------------------------------------------
# Scroll the sheet down a bit and move the active selection the same way
# the author's Excel session left it after adding the row.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 44
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C50").Select() | Out-Null
